{"js": "// Update the 25 \"three-digit \u00d7 one-digit\" equations in the worksheet's\n// single 20-row table (data lives on rows 0, 4, 9, 14, 19 \u2014 five\n// equations per row, the intervening rows are blank answer rows).\n//\n// We address cells by (row, col) rather than by text search because a\n// couple of the original equations are duplicated verbatim (e.g.\n// \"310\u00d77=2170\" appears twice) but map to DIFFERENT replacement values,\n// so positional addressing is the only unambiguous strategy.\n\nconst dataRowIndexes = [0, 4, 9, 14, 19];\n\n// Each inner array holds the 5 new equations (left-to-right) for the\n// corresponding row in dataRowIndexes, in document order.\nconst replacements = [\n  [\"599\u00d73=1797\", \"613\u00d74=2452\", \"307\u00d73=921\", \"771\u00d78=6168\", \"762\u00d74=3048\"],\n  [\"317\u00d79=2853\", \"119\u00d77=833\", \"537\u00d77=3759\", \"887\u00d72=1774\", \"453\u00d73=1359\"],\n  [\"638\u00d77=4466\", \"749\u00d75=3745\", \"263\u00d77=1841\", \"255\u00d72=510\", \"632\u00d79=5688\"],\n  [\"295\u00d73=885\", \"839\u00d78=6712\", \"859\u00d75=4295\", \"483\u00d79=4347\", \"723\u00d79=6507\"],\n  [\"850\u00d73=2550\", \"803\u00d73=2409\", \"876\u00d73=2628\", \"954\u00d74=3816\", \"578\u00d72=1156\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let i = 0; i < dataRowIndexes.length; i++) {\n  const rowIndex = dataRowIndexes[i];\n  const rowValues = replacements[i];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(rowIndex, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 \"three-digit x one-digit\" equations in the worksheet's\n# single 20-row table (data lives on rows 1, 5, 10, 15, 20 -- 1-indexed,\n# five equations per row; the intervening rows are blank answer rows).\n#\n# Cells are addressed by (row, column) rather than by text search because\n# a couple of the original equations are duplicated verbatim (e.g.\n# \"310x7=2170\" appears twice) but map to DIFFERENT replacement values, so\n# positional addressing is the only unambiguous strategy.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$dataRows = @(1, 5, 10, 15, 20)\n\n$replacements = @(\n    @(\"599\u00d73=1797\", \"613\u00d74=2452\", \"307\u00d73=921\", \"771\u00d78=6168\", \"762\u00d74=3048\"),\n    @(\"317\u00d79=2853\", \"119\u00d77=833\", \"537\u00d77=3759\", \"887\u00d72=1774\", \"453\u00d73=1359\"),\n    @(\"638\u00d77=4466\", \"749\u00d75=3745\", \"263\u00d77=1841\", \"255\u00d72=510\", \"632\u00d79=5688\"),\n    @(\"295\u00d73=885\", \"839\u00d78=6712\", \"859\u00d75=4295\", \"483\u00d79=4347\", \"723\u00d79=6507\"),\n    @(\"850\u00d73=2550\", \"803\u00d73=2409\", \"876\u00d73=2628\", \"954\u00d74=3816\", \"578\u00d72=1156\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n    $rowIndex = $dataRows[$i]\n    $rowValues = $replacements[$i]\n    for ($col = 1; $col -le $rowValues.Length; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
